# Atualização de bases das ligas, do dia: 06-04-2024 às 01:36
# Updates odds/result data for "Romania Liga I" sheet rows 235-239 and 259-264.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 235 (id 233) ---
$ws.Range('B235').Value = 6870268
$ws.Range('F235').Value = 'Petrolul Ploiesti'
$ws.Range('G235').Value = 'ACS Sepsi'
$ws.Range('H235').Value = 1
$ws.Range('I235').Value = 2
$ws.Range('J235').Value = 'A'
$ws.Range('K235').Value = 2.8
$ws.Range('L235').Value = 3
$ws.Range('M235').Value = 2.55
$ws.Range('N235').Value = 3
$ws.Range('O235').Value = 3.2
$ws.Range('P235').Value = 2.3
$ws.Range('R235').Value = 1.85
$ws.Range('S235').Value = 2
$ws.Range('T235').Value = 2.25
$ws.Range('U235').Value = 1.875
$ws.Range('V235').Value = 1.975
$ws.Range('X235').Value = -1
$ws.Range('Y235').Value = 1.3
$ws.Range('Z235').Value = -1
$ws.Range('AA235').Value = 1
$ws.Range('AB235').Value = 0.875
$ws.Range('AC235').Value = -1

# --- Row 236 (id 234) ---
$ws.Range('B236').Value = 6861095
$ws.Range('F236').Value = 'FC Botosani'
$ws.Range('G236').Value = 'Farul Constanta'
$ws.Range('K236').Value = 3.75
$ws.Range('L236').Value = 3.4
$ws.Range('M236').Value = 1.909
$ws.Range('N236').Value = 3.1
$ws.Range('O236').Value = 3
$ws.Range('P236').Value = 2.375
$ws.Range('R236').Value = 1.775
$ws.Range('S236').Value = 2.1
$ws.Range('T236').Value = 2
$ws.Range('U236').Value = 1.8
$ws.Range('V236').Value = 2.05
$ws.Range('X236').Value = 2
$ws.Range('Z236').Value = 0.3875
$ws.Range('AC236').Value = 1.05

# --- Row 237 (id 235) ---
$ws.Range('B237').Value = 6865915
$ws.Range('F237').Value = 'FC Voluntari'
$ws.Range('G237').Value = 'Universitatea Cluj'
$ws.Range('H237').Value = 0
$ws.Range('J237').Value = 'D'
$ws.Range('K237').Value = 3.5
$ws.Range('L237').Value = 3.25
$ws.Range('M237').Value = 2.05
$ws.Range('N237').Value = 3.4
$ws.Range('O237').Value = 3.1
$ws.Range('P237').Value = 2.15
$ws.Range('Q237').Value = 0.25
$ws.Range('R237').Value = 1.975
$ws.Range('S237').Value = 1.875
$ws.Range('U237').Value = 2.05
$ws.Range('V237').Value = 1.75
$ws.Range('W237').Value = -1
$ws.Range('X237').Value = 2.1
$ws.Range('Z237').Value = 0.4875
$ws.Range('AC237').Value = 0.75

# --- Row 238 (id 236) ---
$ws.Range('B238').Value = 6836277
$ws.Range('F238').Value = 'CFR Cluj'
$ws.Range('G238').Value = 'AFC Hermannstadt'
$ws.Range('K238').Value = 1.7
$ws.Range('L238').Value = 3.4
$ws.Range('M238').Value = 5
$ws.Range('N238').Value = 1.65
$ws.Range('O238').Value = 3.5
$ws.Range('P238').Value = 5.25
$ws.Range('Q238').Value = -0.75
$ws.Range('R238').Value = 1.85
$ws.Range('S238').Value = 2
$ws.Range('U238').Value = 1.875
$ws.Range('V238').Value = 1.975
$ws.Range('W238').Value = 0.6499999999999999
$ws.Range('Z238').Value = 0.425
$ws.Range('AA238').Value = -0.5
$ws.Range('AC238').Value = 0.9750000000000001

# --- Row 239 (id 237) ---
$ws.Range('B239').Value = 6852370
$ws.Range('F239').Value = 'Dinamo Bucharest'
$ws.Range('G239').Value = 'ACS UTA Batrana Doamna'
$ws.Range('I239').Value = 0
$ws.Range('J239').Value = 'H'
$ws.Range('K239').Value = 2.55
$ws.Range('L239').Value = 2.875
$ws.Range('M239').Value = 3
$ws.Range('N239').Value = 2.375
$ws.Range('O239').Value = 3
$ws.Range('P239').Value = 3.1
$ws.Range('Q239').Value = -0.25
$ws.Range('R239').Value = 2
$ws.Range('S239').Value = 1.85
$ws.Range('U239').Value = 1.975
$ws.Range('V239').Value = 1.875
$ws.Range('W239').Value = 1.375
$ws.Range('Y239').Value = -1
$ws.Range('Z239').Value = 1
$ws.Range('AA239').Value = -1
$ws.Range('AB239').Value = -1
$ws.Range('AC239').Value = 0.875

# --- Row 259 (id 257) : match result now known ---
$ws.Range('H259').Value = 1
$ws.Range('I259').Value = 0
$ws.Range('J259').Value = 'H'
$ws.Range('N259').Value = 2.2
$ws.Range('P259').Value = 3.3
$ws.Range('Q259').Value = -0.25
$ws.Range('R259').Value = 1.925
$ws.Range('S259').Value = 1.925
$ws.Range('T259').Value = 2.5
$ws.Range('U259').Value = 2
$ws.Range('V259').Value = 1.85
$ws.Range('W259').Value = 1.2
$ws.Range('X259').Value = -1
$ws.Range('Y259').Value = -1
$ws.Range('Z259').Value = 0.925
$ws.Range('AA259').Value = -1
$ws.Range('AB259').Value = -1
$ws.Range('AC259').Value = 0.8500000000000001

# --- Row 260 (id 258) : match result now known ---
$ws.Range('H260').Value = 1
$ws.Range('I260').Value = 1
$ws.Range('J260').Value = 'D'
$ws.Range('N260').Value = 2
$ws.Range('O260').Value = 3.3
$ws.Range('P260').Value = 3.6
$ws.Range('Q260').Value = -0.5
$ws.Range('R260').Value = 2.05
$ws.Range('S260').Value = 1.8
$ws.Range('W260').Value = -1
$ws.Range('X260').Value = 2.3
$ws.Range('Y260').Value = -1
$ws.Range('Z260').Value = -1
$ws.Range('AA260').Value = 0.8
$ws.Range('AB260').Value = -0.5
$ws.Range('AC260').Value = 0.5125

# --- Row 261 (id 259) : odds movement only ---
$ws.Range('R261').Value = 2.025
$ws.Range('S261').Value = 1.825

# --- Row 262 (id 260) : odds movement only ---
$ws.Range('R262').Value = 2.025
$ws.Range('S262').Value = 1.825
$ws.Range('U262').Value = 1.825
$ws.Range('V262').Value = 2.025

# --- Row 263 (id 261) : odds movement only ---
$ws.Range('U263').Value = 2
$ws.Range('V263').Value = 1.85

# --- Row 264 (id 262) : odds movement only ---
$ws.Range('R264').Value = 1.975
$ws.Range('S264').Value = 1.875
